$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2..25, columns B, C, D, E, F, I, L
# (columns G, H, J, K, M, N, O remain unchanged at 0)
$newValues = @(
    @(1.792708522336056, 0.2544466942002543, 0.08884395901083764, 0.04573766239177823, 2.122403679451622, 1.546942915480017, 0.2749935328723581),
    @(1.678671900351333, 0.2219482845705159, 0.08912226837344761, 0.04589353640667726, 2.068054950713929, 1.525602533049664, 0.2647931738955265),
    @(1.609751299265668, 0.2020316225884358, 0.08932779141339253, 0.04599787337157046, 2.036019618066291, 1.513369054356119, 0.2587060589918622),
    @(1.581940629925157, 0.1939241500997468, 0.08942018263914875, 0.04604256283601282, 2.023297626148391, 1.508600719692922, 0.2562695208956995),
    @(1.577339268476578, 0.1925784175072636, 0.08943604381180137, 0.04605011466809605, 2.021205153028461, 1.50782199065921, 0.2558675870444063),
    @(1.609375122740403, 0.2019222479428322, 0.08932900254993115, 0.0459984672758309, 2.035846701905967, 1.51330387102584, 0.2586730210912407),
    @(1.753159983072294, 0.2432330273924208, 0.08893268737445581, 0.04578961742616716, 2.103385306817557, 1.539403194276687, 0.2714397861478943),
    @(2.04390787506776, 0.3245766972721071, 0.08843361450761478, 0.0454484832348605, 2.246566849252986, 1.597567314312627, 0.2978837439293898),
    @(2.263003415741764, 0.3845990521124918, 0.08824100734488383, 0.04523949753718881, 2.358523981160857, 1.644679326439928, 0.3181903824846017),
    @(2.363896293203538, 0.411973697710323, 0.08819216580262434, 0.04515345536416948, 2.410971984475225, 1.667089367148805, 0.3276235987047187),
    @(2.402279996848847, 0.4223507785983429, 0.08817932775077963, 0.0451221706047864, 2.431054634634023, 1.675718259728015, 0.3312241668824072),
    @(2.394005447673919, 0.4201153909158393, 0.08818183974981508, 0.04512885063549099, 2.426719560200922, 1.673853495084188, 0.3304474518009783),
    @(2.367050571730317, 0.4128272033377698, 0.08819099576307821, 0.04515085554580534, 2.412619732660147, 1.667796400064091, 0.3279192485103408),
    @(2.350563129600289, 0.4083644212853983, 0.0881973432373897, 0.04516450315515197, 2.404012159129366, 1.664104896915021, 0.3263743598546256),
    @(2.256434616480419, 0.3828115352195596, 0.08824498645797973, 0.04524530219111078, 2.35512722140902, 1.643234637527456, 0.3175778601507204),
    @(2.19900472838242, 0.3671542086020736, 0.08828420313628271, 0.04529718133626748, 2.325529151467606, 1.630683444075203, 0.3122317960713161),
    @(2.166087800250295, 0.3581550844531876, 0.08831040445684124, 0.0453278704854867, 2.308647812362238, 1.623556315742576, 0.3091752771527467),
    @(2.154962437060135, 0.355109237383715, 0.08831989933112538, 0.04533840724496141, 2.30295648409907, 1.621158931153388, 0.308143543330047),
    @(2.205106301524779, 0.3688202734474544, 0.0882796506320247, 0.04529157078200896, 2.328665124218873, 1.632010005528187, 0.3127989877993116),
    @(2.374963034898826, 0.41496761845832, 0.08818815223504117, 0.04514435695853969, 2.416755148649258, 1.669571627063178, 0.3286610696633829),
    @(2.487011137445734, 0.4451914768026768, 0.08816136531720531, 0.04505570713376894, 2.475621097480911, 1.694952996079621, 0.3391936116696002),
    @(2.42711353955076, 0.4290543171631498, 0.08817261433283363, 0.04510232933776681, 2.444083688409336, 1.681329631784152, 0.3335569337935453),
    @(2.202347468076084, 0.3680670378622608, 0.0882816974406424, 0.0452941046251345, 2.327246930746355, 1.631409990902043, 0.312542507475186),
    @(1.964299393444946, 0.3025296680355325, 0.08853838897886135, 0.04553345119753471, 2.206660695328992, 1.581072248608024, 0.2905771658949305)
)

$columns = @(2, 3, 4, 5, 6, 9, 12)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $startRow + $i
    $rowValues = $newValues[$i]
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $col = $columns[$j]
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

Write-Output "Updated $($newValues.Count) rows"
